$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 93.419141172176538
$ws.Range("C2").Value = 93.321681545976858
$ws.Range("D2").Value = 94.292850339544714
$ws.Range("E2").Value = 94.550870470827235

$ws.Range("B3").Value = 93.786045175059073
$ws.Range("C3").Value = 97.147975364030302
$ws.Range("D3").Value = 95.751039117096752
$ws.Range("E3").Value = 94.57436246392929

$ws.Range("B1:E3").Select()
